$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the QUANTITY column values to 100
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 100
$ws.Range("C4").Value = 100

# Update the HEADER CODE for row 3 from "IA 01" to "IA 02"
$ws.Range("A3").Value = "IA 02"

# Move the active selection to C5
$ws.Range("C5").Select()
